# "Secondaire => Niveau partout"
#
# The running footer abbreviates the section/level label as
# "pyETUDE_Sec" (a merge-field-style placeholder meaning "Secondaire").
# Rename it to "pyETUDE_Niv" ("Niveau") everywhere it shows up in the
# document -- body text, and every header/footer of every section
# (primary, first-page and even-page variants), so the rename is
# applied "partout" regardless of which story happens to carry it.

$d = $word.ActiveDocument

function Replace-Everywhere($range, [string]$find, [string]$replace) {
    if ($range -eq $null) { return $false }
    return $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $replace, 2)
}

$oldText = "pyETUDE_Sec"
$newText = "pyETUDE_Niv"

# Main document story (title page, table of contents, body, ...).
Replace-Everywhere $d.Content $oldText $newText | Out-Null

# Every section's headers/footers -- index 1/2/3 map to the primary,
# first-page and even-page variants respectively.
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    for ($idx = 1; $idx -le 3; $idx++) {
        $footer = $section.Footers($idx)
        if ($footer -ne $null -and $footer.Exists) {
            Replace-Everywhere $footer.Range $oldText $newText | Out-Null
        }

        $header = $section.Headers($idx)
        if ($header -ne $null -and $header.Exists) {
            Replace-Everywhere $header.Range $oldText $newText | Out-Null
        }
    }
}

# Footnotes / endnotes / comments, in case the placeholder ever ends up
# there too (no-ops on this document, but keeps the rename thorough).
for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    Replace-Everywhere $d.Footnotes($i).Range $oldText $newText | Out-Null
}
for ($i = 1; $i -le $d.Endnotes.Count; $i++) {
    Replace-Everywhere $d.Endnotes($i).Range $oldText $newText | Out-Null
}
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    Replace-Everywhere $d.Comments($i).Range $oldText $newText | Out-Null
}
